# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (7-column fund holdings detail) right
#    before the "总计" (totals) sheet, duplicating "总计"'s formatting so
#    the header/index-column style (bold, thin border, centered) matches
#    the other detail sheets.
# 2. Insert a new first data row into "总计" summarising the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Helper: force a cell to store its value as literal text (keeps leading
# zeros / avoids Excel's automatic number coercion for things like fund
# codes "006679" or text-formatted ratios "14.75"), then drop the
# NumberFormat-driven style residue so the cell ends up with no explicit
# style - matching the look of the other plain data cells in this
# workbook.
# ---------------------------------------------------------------------
function Set-TextValue($rng, $val) {
  $rng.NumberFormat = "@"
  $rng.Value = $val
  $rng.ClearFormats()
}

# -----------------------------------------------------------------
# Step 1: duplicate "总计" to get a correctly-styled new sheet, placed
# immediately before "总计" itself, then rename + rebuild its contents.
# -----------------------------------------------------------------
$total.Copy($total)
# Re-resolve "总计" - the old $total handle's .Index is stale after the
# sheet collection shifted, so look it up by name again.
$total = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Item($total.Index - 1)
$ws.Name = "2022-Q1"

$ws.Cells.ClearContents()

# Extend the header style (already s=2 from the copied sheet on B1:D1)
# across the extra columns E:H, and extend the index-column style
# (already s=2 on A2:A4) down to the extra row A5.
$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Row 2
$ws.Range("A2").Value = 0
Set-TextValue $ws.Range("B2") "006679"
$ws.Range("C2").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A"
Set-TextValue $ws.Range("D2") "14.75"
Set-TextValue $ws.Range("E2") "83.19"
Set-TextValue $ws.Range("F2") "3.54"
Set-TextValue $ws.Range("G2") "0.5222"
$ws.Range("H2").Value = 8

# Row 3
$ws.Range("A3").Value = 1
Set-TextValue $ws.Range("B3") "162719"
$ws.Range("C3").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
Set-TextValue $ws.Range("D3") "14.75"
Set-TextValue $ws.Range("E3") "83.19"
Set-TextValue $ws.Range("F3") "3.54"
Set-TextValue $ws.Range("G3") "0.5222"
$ws.Range("H3").Value = 8

# Row 4
$ws.Range("A4").Value = 2
Set-TextValue $ws.Range("B4") "006680"
$ws.Range("C4").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C"
Set-TextValue $ws.Range("D4") "4.73"
Set-TextValue $ws.Range("E4") "83.19"
Set-TextValue $ws.Range("F4") "3.54"
Set-TextValue $ws.Range("G4") "0.1674"
$ws.Range("H4").Value = 8

# Row 5
$ws.Range("A5").Value = 3
Set-TextValue $ws.Range("B5") "004243"
$ws.Range("C5").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
Set-TextValue $ws.Range("D5") "4.73"
Set-TextValue $ws.Range("E5") "83.19"
Set-TextValue $ws.Range("F5") "3.54"
Set-TextValue $ws.Range("G5") "0.1674"
$ws.Range("H5").Value = 8

# -----------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计"'s data
# (row 2), pushing the existing quarters down, and renumber the index
# column (A) so it stays 0,1,2,3.
# -----------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 1.38

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the originally active sheet/tab (Sheet.Copy() above shifts the
# active tab onto the freshly duplicated sheet as a side effect).
$wb.Worksheets.Item(1).Activate()

Write-Host "2022-Q1 sheet added and 总计 updated"
